$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 252, shifting rows 252:273 down to 253:274
$ws.Rows.Item(252).Insert()

# Populate the newly inserted row 252 with the new weekly data point.
$ws.Cells.Item(252, 1).Value = 8
$ws.Cells.Item(252, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(252, 3).Value = "Coquimbo"
$ws.Cells.Item(252, 4).Value = 45265
$ws.Cells.Item(252, 5).Value = 4
$ws.Cells.Item(252, 6).Value = 100112040
$ws.Cells.Item(252, 7).Value = "Cilantro"
$ws.Cells.Item(252, 8).Value = "Sin especificar"
$ws.Cells.Item(252, 9).Value = "Primera"
$ws.Cells.Item(252, 10).Value = 3000
$ws.Cells.Item(252, 11).Value = 2400
$ws.Cells.Item(252, 12).Value = 2500
$ws.Cells.Item(252, 13).Value = 2450
$ws.Cells.Item(252, 14).Value = "$/atado 1 a 1,5 kilos"
$ws.Cells.Item(252, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(252, 16).Value = 1633
$ws.Cells.Item(252, 17).Value = 1.5
$ws.Cells.Item(252, 18).Value = "Hortaliza"
